$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.479.62'
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.571.97'
$ws.Range("E3").Value = '  +0.08%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("E5").Value = '  -0.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '291.95'
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3723'
$ws.Range("E7").Value = '  -0.99%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.95'
$ws.Range("E8").Value = '  +0.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3400'
$ws.Range("E9").Value = '  -0.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.144'
$ws.Range("E10").Value = '  +0.13%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07553'
$ws.Range("E11").Value = '  -0.83%  '
$ws.Range("E12").Value = '  -0.07%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.30'
$ws.Range("E13").Value = '  +0.79%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.041'
$ws.Range("E14").Value = '  +0.57%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.969'
$ws.Range("E15").Value = '  +0.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.569.95'
$ws.Range("E16").Value = '  -0.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001123'
$ws.Range("E17").Value = '  -0.94%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.74'
$ws.Range("E18").Value = '  +0.65%  '
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.293'
$ws.Range("E21").Value = '  +1.63%  '
$ws.Range("E22").Value = '  -2.20%  '
$ws.Range("E23").Value = '  +1.78%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.483.20'
$ws.Range("E24").Value = '  +0.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.371'
$ws.Range("E25").Value = '  -0.76%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.641'
$ws.Range("E26").Value = '  -1.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.02'
$ws.Range("E27").Value = '  -0.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '149.19'
$ws.Range("E28").Value = '  +1.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.063'
$ws.Range("E29").Value = '  +0.62%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.18'
$ws.Range("E30").Value = '  -1.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.747.33'
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.089'
$ws.Range("E32").Value = '  +8.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.242'
$ws.Range("E33").Value = '  +2.51%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.009'
$ws.Range("E34").Value = '  -0.38%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.783'
$ws.Range("E35").Value = '  -3.46%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08377'
$ws.Range("E36").Value = '  -1.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02489'
$ws.Range("E37").Value = '  -1.99%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2305'
$ws.Range("E38").Value = '  +0.01%  '
$ws.Range("E39").Value = '  -3.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06525'
$ws.Range("E40").Value = '  +0.30%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.458'
$ws.Range("E41").Value = '  +0.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.39'
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6237'
$ws.Range("E43").Value = '  -1.49%  '
$ws.Range("E44").Value = '  -0.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.03'
$ws.Range("E45").Value = '  +0.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.813'
$ws.Range("E46").Value = '  +0.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5868'
$ws.Range("E47").Value = '  -1.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '130.70'
$ws.Range("E48").Value = '  +5.05%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.074'
$ws.Range("E49").Value = '  -0.49%  '
$ws.Range("E50").Value = '  -5.16%  '
$ws.Range("E51").Value = '  +0.19%  '
